# Adds two new columns ("Осталось дней" / "Примечание") with certificate
# expiry/error information to the two "report" worksheets.
#
# Sheet "report.xml" (index 1) and "report — копия.xml" (index 2) both
# gain columns F and G for every row. The third sheet ("Sheet") is left
# untouched.

$wb = $excel.ActiveWorkbook

$sheetIndexes = 1, 2

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # --- Header row ---------------------------------------------------
    $ws.Range("F1").Value = "Осталось дней"
    $ws.Range("G1").Value = "Примечание"

    # --- Row 2 (87.240.190.72 / Sectigo) -------------------------------
    $ws.Range("F2").Value = 121
    $ws.Range("G2").Value = "<X509Name object '/C=GB/ST=Greater Manchester/L=Salford/O=Sectigo Limited/CN=Sectigo ECC Extended Validation Secure Server CA'>"

    # --- Column widths (F=22 chars, G=100 chars) -----------------------
    # ColumnWidth snaps to the host's pixel grid (MDW=7), so the raw
    # values below are chosen to round-trip to exactly 22 / 100.
    $ws.Columns.Item(6).ColumnWidth = 21.16
    $ws.Columns.Item(7).ColumnWidth = 99.19
}

# --- Sheet "report.xml": row 3 (5.255.255.80 / Yandex CA) --------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 178
$ws1.Range("G3").Value = "<X509Name object '/C=RU/O=Yandex LLC/OU=Yandex Certification Authority/CN=Yandex CA'>"

# --- Sheet "report — копия.xml": row 3 (10.1.2.249 / unreachable) ------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C3").Value = ""
$ws2.Range("D3").Value = ""
$ws2.Range("E3").Value = ""
$ws2.Range("F3").Value = ""
$ws2.Range("G3").Value = "[WinError 10060] Попытка установить соединение была безуспешной, т.к. от другого компьютера за требуемое время не получен нужный отклик, или было разорвано уже установленное соединение из-за неверного отклика уже подключенного компьютера"

Write-Output "edit applied"
